$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dates  = @(44326, 44327, 44328, 44329)
$newPos = @(0, 0, 0, 1)
$somma  = @(8, 8, 8, 6)
$per100 = @(244.2748091603054, 244.2748091603054, 244.2748091603054, 183.206106870229)

$startRow = 252
for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $startRow + $i
    $prev = $r - 1

    # Copy formatting from the row above (keeps existing styles, e.g. date format on col A)
    $ws.Cells.Item($prev, 1).Copy($ws.Cells.Item($r, 1))

    # Now overwrite the values for this new row
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $newPos[$i]
    $ws.Cells.Item($r, 3).Value = $somma[$i]
    $ws.Cells.Item($r, 4).Value = $per100[$i]
}
